$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K7").Value = 0.1683237681281231
$ws.Range("J8").Value = 0.1722916656412322
$ws.Range("I9").Value = 0.3079317558114735
$ws.Range("H10").Value = 0.06712557395580883
$ws.Range("G11").Value = 0.02179435870371246
$ws.Range("F12").Value = -0.04506706323234141
$ws.Range("E13").Value = -0.07465326558905801
$ws.Range("D14").Value = -0.0928039223186989
$ws.Range("C15").Value = -0.1108357465673982
$ws.Range("B16").Value = -0.1624199859130616
